$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn 1 bác sĩ")

# Insert a new "Nhóm dịch vụ" column before the existing "Tên dịch vụ" column (G),
# shifting every following column one place to the right.
$ws.Columns("G").Insert()
$ws.Range("G1").Value = "Nhóm dịch vụ"
$ws.Range("G2").Value = "Ngực"

# Row 3 (the totals row) had no text in this column originally - keep it a truly
# blank cell (matching the other blank cells on that row) instead of leaving it
# with the "0" that a freshly inserted column would otherwise carry.
$ws.Range("G3").ClearContents()
$ws.Range("Z99").Copy($ws.Range("G3"))

# Add the four new trailing columns used for the sale commission breakdown.
$ws.Range("X1").Value = "Tỉ lệ chiết khấu sale chính"
$ws.Range("Y1").Value = "Tỉ lệ chiết khấu sale phụ"
$ws.Range("Z1").Value = "Chiết khấu sale chính"
$ws.Range("AA1").Value = "Chiết khấu sale phụ"

$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0

$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 0
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 0

# The payroll report is no longer generated from this workbook - drop the sheet.
$wb.Worksheets.Item("Lương").Delete()
